$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, pushing existing rows 46..83 down to 47..84
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new weekly data point
$ws.Range("A46").Value = 9
$ws.Range("B46").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C46").Value = "Metropolitana"
$ws.Range("D46").Value = 45072
$ws.Range("E46").Value = 13
$ws.Range("F46").Value = 100112035
$ws.Range("G46").Value = "Bruselas (repollito)"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 43
$ws.Range("K46").Value = 22000
$ws.Range("L46").Value = 24000
$ws.Range("M46").Value = 23023
$ws.Range("N46").Value = "$/malla 15 kilos"
$ws.Range("O46").Value = "Provincia de Quillota"
$ws.Range("P46").Value = 1535
$ws.Range("Q46").Value = 15
$ws.Range("R46").Value = "Hortaliza"
